$wb = $excel.ActiveWorkbook

# --- frequent_itemsets1 sheet: rows 3-50 reordered ---
$ws1 = $wb.Worksheets.Item("frequent_itemsets1")
$ws1.Range("B3").Value = "{9}"
$ws1.Range("C3").Value = 225
$ws1.Range("B4").Value = "{15}"
$ws1.Range("C4").Value = 185
$ws1.Range("B5").Value = "{47}"
$ws1.Range("C5").Value = 34
$ws1.Range("B6").Value = "{36}"
$ws1.Range("C6").Value = 174
$ws1.Range("B7").Value = "{10}"
$ws1.Range("C7").Value = 34
$ws1.Range("B8").Value = "{46}"
$ws1.Range("C8").Value = 48
$ws1.Range("B9").Value = "{45}"
$ws1.Range("C9").Value = 32
$ws1.Range("B10").Value = "{32}"
$ws1.Range("C10").Value = 36
$ws1.Range("B11").Value = "{5}"
$ws1.Range("C11").Value = 35
$ws1.Range("B12").Value = "{40}"
$ws1.Range("C12").Value = 42
$ws1.Range("B13").Value = "{41}"
$ws1.Range("C13").Value = 51
$ws1.Range("B14").Value = "{43}"
$ws1.Range("C14").Value = 46
$ws1.Range("B15").Value = "{42}"
$ws1.Range("C15").Value = 36
$ws1.Range("B16").Value = "{21}"
$ws1.Range("C16").Value = 49
$ws1.Range("B17").Value = "{14}"
$ws1.Range("C17").Value = 293
$ws1.Range("B18").Value = "{35}"
$ws1.Range("C18").Value = 43
$ws1.Range("B19").Value = "{1}"
$ws1.Range("C19").Value = 156
$ws1.Range("B20").Value = "{37}"
$ws1.Range("C20").Value = 37
$ws1.Range("B21").Value = "{11}"
$ws1.Range("C21").Value = 44
$ws1.Range("B22").Value = "{34}"
$ws1.Range("C22").Value = 33
$ws1.Range("B23").Value = "{8}"
$ws1.Range("C23").Value = 41
$ws1.Range("B24").Value = "{18}"
$ws1.Range("C24").Value = 42
$ws1.Range("B25").Value = "{12}"
$ws1.Range("C25").Value = 300
$ws1.Range("B26").Value = "{26}"
$ws1.Range("C26").Value = 47
$ws1.Range("B27").Value = "{31}"
$ws1.Range("C27").Value = 40
$ws1.Range("B28").Value = "{0}"
$ws1.Range("C28").Value = 32
$ws1.Range("B29").Value = "{29}"
$ws1.Range("C29").Value = 36
$ws1.Range("B30").Value = "{44}"
$ws1.Range("C30").Value = 38
$ws1.Range("B31").Value = "{49}"
$ws1.Range("C31").Value = 162
$ws1.Range("B32").Value = "{20}"
$ws1.Range("C32").Value = 39
$ws1.Range("B33").Value = "{7}"
$ws1.Range("C33").Value = 42
$ws1.Range("B34").Value = "{30}"
$ws1.Range("C34").Value = 37
$ws1.Range("B35").Value = "{48}"
$ws1.Range("C35").Value = 44
$ws1.Range("B36").Value = "{22}"
$ws1.Range("C36").Value = 215
$ws1.Range("B37").Value = "{27}"
$ws1.Range("C37").Value = 38
$ws1.Range("B38").Value = "{13}"
$ws1.Range("C38").Value = 36
$ws1.Range("B39").Value = "{3}"
$ws1.Range("C39").Value = 34
$ws1.Range("B40").Value = "{16}"
$ws1.Range("C40").Value = 283
$ws1.Range("B41").Value = "{33}"
$ws1.Range("C41").Value = 47
$ws1.Range("B42").Value = "{2}"
$ws1.Range("C42").Value = 43
$ws1.Range("B43").Value = "{23}"
$ws1.Range("C43").Value = 41
$ws1.Range("B44").Value = "{24}"
$ws1.Range("C44").Value = 37
$ws1.Range("B45").Value = "{6}"
$ws1.Range("C45").Value = 31
$ws1.Range("B46").Value = "{39}"
$ws1.Range("C46").Value = 32
$ws1.Range("B47").Value = "{17}"
$ws1.Range("C47").Value = 34
$ws1.Range("B48").Value = "{25}"
$ws1.Range("C48").Value = 38
$ws1.Range("B49").Value = "{38}"
$ws1.Range("C49").Value = 38
$ws1.Range("B50").Value = "{28}"
$ws1.Range("C50").Value = 34

# --- frequent_itemsets2 sheet: rows 4-8 reordered ---
$ws2 = $wb.Worksheets.Item("frequent_itemsets2")
$ws2.Range("B4").Value = "{15, 36}"
$ws2.Range("C4").Value = 139
$ws2.Range("B5").Value = "{14, 12}"
$ws2.Range("C5").Value = 268
$ws2.Range("B6").Value = "{14, 16}"
$ws2.Range("C6").Value = 258
$ws2.Range("B7").Value = "{1, 49}"
$ws2.Range("C7").Value = 127
$ws2.Range("B8").Value = "{12, 16}"
$ws2.Range("C8").Value = 259

# --- association_rules sheet: rows 5-14 reordered ---
$ws4 = $wb.Worksheets.Item("association_rules")
$ws4.Range("B5").Value = "{15}"
$ws4.Range("C5").Value = "{36}"
$ws4.Range("D5").Value = 0.7513513513513513
$ws4.Range("B6").Value = "{36}"
$ws4.Range("C6").Value = "{15}"
$ws4.Range("D6").Value = 0.7988505747126436
$ws4.Range("B7").Value = "{14}"
$ws4.Range("C7").Value = "{12}"
$ws4.Range("D7").Value = 0.9146757679180887
$ws4.Range("B8").Value = "{12}"
$ws4.Range("C8").Value = "{14}"
$ws4.Range("D8").Value = 0.8933333333333333
$ws4.Range("B9").Value = "{14}"
$ws4.Range("C9").Value = "{16}"
$ws4.Range("D9").Value = 0.8805460750853242
$ws4.Range("B10").Value = "{16}"
$ws4.Range("C10").Value = "{14}"
$ws4.Range("D10").Value = 0.911660777385159
$ws4.Range("B11").Value = "{1}"
$ws4.Range("C11").Value = "{49}"
$ws4.Range("D11").Value = 0.8141025641025641
$ws4.Range("B12").Value = "{49}"
$ws4.Range("C12").Value = "{1}"
$ws4.Range("D12").Value = 0.7839506172839507
$ws4.Range("B13").Value = "{12}"
$ws4.Range("C13").Value = "{16}"
$ws4.Range("D13").Value = 0.8633333333333333
$ws4.Range("B14").Value = "{16}"
$ws4.Range("C14").Value = "{12}"
$ws4.Range("D14").Value = 0.9151943462897526
